$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 36204183
$ws.Range("B2").Value = 16015905
$ws.Range("C2").Value = 50522577
$ws.Range("D2").Value = 26543024
$ws.Range("E2").Value = 14318394
$ws.Range("F2").Value = 39.55
$ws.Range("G2").Value = 10527119
$ws.Range("H2").Value = 65.73
